# Actualización automática 2025-10-20 15:30:09
# Updates sales/compliance figures for CASTRO ALCIVAR EDA MARIA across the
# three report sheets (VENTAS POR GRUPO, VENTA MENSUAL, CUMPLIMIENTO MENSUAL).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 33 (CASTRO ALCIVAR EDA MARIA / MAD&DECO S.A.)
$wsGrupo.Range("D33").Value = 1780.8
$wsGrupo.Range("K33").Value = 1006.4
$wsGrupo.Range("M33").Value = 897.4400000000001

# Row 48 (CASTRO ALCIVAR EDA MARIA / SANCHEZ CORREA MARCO EDUARDO)
$wsGrupo.Range("M48").Value = 272.02

# Row 59 (totals "N de 57")
$wsGrupo.Range("K59").Value = "1 de 57"
$wsGrupo.Range("M59").Value = "7 de 57"

# --- Sheet 2: VENTA MENSUAL --------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F33").Value = 3684.64
$wsMensual.Range("F48").Value = 804.95
$wsMensual.Range("F59").Value = 22339.39

# --- Sheet 3: CUMPLIMIENTO MENSUAL ------------------------------------------
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3 (240X80 PORCELANATO)
$wsCumpl.Range("D3").Value = 4147.59
$wsCumpl.Range("E3").Value = 16239.8874217135
$wsCumpl.Range("F3").Value = 0.2034381161635351

# Row 10 (PANELES DECORATIVOS)
$wsCumpl.Range("D10").Value = 1062.72
$wsCumpl.Range("E10").Value = 1654.03588474074
$wsCumpl.Range("F10").Value = 0.3911724295763936

# Row 12 (PORCELANATO)
$wsCumpl.Range("D12").Value = 9157.879999999999
$wsCumpl.Range("E12").Value = 39466.18
$wsCumpl.Range("F12").Value = 0.1883405046801933

# Row 14 (TOTAL)
$wsCumpl.Range("D14").Value = 25908.01
$wsCumpl.Range("E14").Value = 73989.98284188786
$wsCumpl.Range("F14").Value = 0.2593446501072904
